$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.254.95"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.787.71"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3786"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3442"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.201"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07512"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.489"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.07%  "

$ws.Range("D15").Value = "1.788.13"
$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.112"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001101"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06668"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.33%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.640"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.44%  "

$ws.Range("D23").Value = "27.247.93"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.507"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.560"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "1.988.73"
$ws.Range("E30").Value = "  -1.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.133"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.661"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6968"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.481"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2211"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.846"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06347"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02341"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.245"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.153"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07152"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
